$wb = $excel.ActiveWorkbook

# --- Provencher sheet: log two new timesheet entries (rows 20-21) ---
$ws2 = $wb.Worksheets.Item("Provencher")

$ws2.Range("A20").Value = 42322
$ws2.Range("B20").Value = 5
$ws2.Range("C20").Value = "Ressearch into raycast assisted movement stabalization"

$ws2.Range("A21").Value = 42323
$ws2.Range("B21").Value = 3.5
$ws2.Range("C21").Value = "Implementation of racast assisted movement stablization"

# --- Team Meetings sheet: log two new meeting entries (rows 17-18) ---
$ws1 = $wb.Worksheets.Item("Team Meetings")

$ws1.Range("A17").Value = 42318
$ws1.Range("B17").Value = 1
$ws1.Range("B17").NumberFormat = "0.00"
$ws1.Range("C17").Value = "Discussion on remaining tasks and issues"

$ws1.Range("A18").Value = 42321
$ws1.Range("B18").Value = 0.5
$ws1.Range("C18").Value = "Discussion of tasks splitting for weekend work"

# --- Update view/selection state to match the latest editing session ---
# Provencher was the active tab before; the new session leaves Team Meetings
# selected with cursor at C18, and Provencher's cursor at B22.
$ws2.Activate() | Out-Null
$ws2.Range("B22").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("C18").Select() | Out-Null
